$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(64, 1).Value = 'Waiver'
$ws.Cells.Item(64, 2).Value = 'MD'
$ws.Cells.Item(64, 3).Value = '1915(c)'
$ws.Cells.Item(64, 4).Value = 'Amendment'
$ws.Cells.Item(64, 5).Value = 'MD-2260.R00.67'
$ws.Cells.Item(64, 6).Value = 'Submitted'
$ws.Cells.Item(64, 7).Value = 'MD-2260.R00.00'

$ws.Cells.Item(65, 1).Value = 'SPA'
$ws.Cells.Item(65, 2).Value = 'MD'
$ws.Cells.Item(65, 3).Value = 'Medicaid SPA'
$ws.Cells.Item(8, 4).Copy($ws.Cells.Item(65, 4))
$ws.Cells.Item(65, 5).Value = 'MD-25-9543'
$ws.Cells.Item(65, 6).Value = 'Approved'
$ws.Cells.Item(8, 7).Copy($ws.Cells.Item(65, 7))

$ws.Cells.Item(66, 1).Value = 'SPA'
$ws.Cells.Item(66, 2).Value = 'MD'
$ws.Cells.Item(66, 3).Value = 'Medicaid SPA'
$ws.Cells.Item(8, 4).Copy($ws.Cells.Item(66, 4))
$ws.Cells.Item(66, 5).Value = 'MD-25-9544'
$ws.Cells.Item(66, 6).Value = 'Under Review'
$ws.Cells.Item(8, 7).Copy($ws.Cells.Item(66, 7))

$ws.Cells.Item(67, 1).Value = 'SPA'
$ws.Cells.Item(67, 2).Value = 'MD'
$ws.Cells.Item(67, 3).Value = 'Medicaid SPA'
$ws.Cells.Item(8, 4).Copy($ws.Cells.Item(67, 4))
$ws.Cells.Item(67, 5).Value = 'MD-25-9545'
$ws.Cells.Item(67, 6).Value = 'Disapproved'
$ws.Cells.Item(8, 7).Copy($ws.Cells.Item(67, 7))

$ws.Cells.Item(68, 1).Value = 'Waiver'
$ws.Cells.Item(68, 2).Value = 'MD'
$ws.Cells.Item(68, 3).Value = '1915(c)'
$ws.Cells.Item(68, 4).Value = 'Amendment'
$ws.Cells.Item(68, 5).Value = 'MD-2260.R00.68'
$ws.Cells.Item(68, 6).Value = 'Pending-Approval'
$ws.Cells.Item(68, 7).Value = 'MD-2260.R00.00'

$ws.Cells.Item(69, 1).Value = 'SPA'
$ws.Cells.Item(69, 2).Value = 'MD'
$ws.Cells.Item(69, 3).Value = 'Medicaid SPA'
$ws.Cells.Item(8, 4).Copy($ws.Cells.Item(69, 4))
$ws.Cells.Item(69, 5).Value = 'MD-25-9546'
$ws.Cells.Item(69, 6).Value = 'Pending-Concurrence'
$ws.Cells.Item(8, 7).Copy($ws.Cells.Item(69, 7))

$ws.Cells.Item(70, 1).Value = 'SPA'
$ws.Cells.Item(70, 2).Value = 'MD'
$ws.Cells.Item(70, 3).Value = 'CHIP SPA'
$ws.Cells.Item(8, 4).Copy($ws.Cells.Item(70, 4))
$ws.Cells.Item(70, 5).Value = 'MD-25-9547'
$ws.Cells.Item(70, 6).Value = 'Submitted'
$ws.Cells.Item(8, 7).Copy($ws.Cells.Item(70, 7))

$ws.Cells.Item(71, 1).Value = 'SPA'
$ws.Cells.Item(71, 2).Value = 'MD'
$ws.Cells.Item(71, 3).Value = 'Medicaid SPA'
$ws.Cells.Item(8, 4).Copy($ws.Cells.Item(71, 4))
$ws.Cells.Item(71, 5).Value = 'MD-25-9548'
$ws.Cells.Item(71, 6).Value = 'RAI Issued'
$ws.Cells.Item(8, 7).Copy($ws.Cells.Item(71, 7))

$ws.Cells.Item(72, 1).Value = 'SPA'
$ws.Cells.Item(72, 2).Value = 'MD'
$ws.Cells.Item(72, 3).Value = 'CHIP SPA'
$ws.Cells.Item(8, 4).Copy($ws.Cells.Item(72, 4))
$ws.Cells.Item(72, 5).Value = 'MD-25-9549'
$ws.Cells.Item(72, 6).Value = 'Submitted'
$ws.Cells.Item(8, 7).Copy($ws.Cells.Item(72, 7))

$ws.Cells.Item(73, 1).Value = 'SPA'
$ws.Cells.Item(73, 2).Value = 'MD'
$ws.Cells.Item(73, 3).Value = 'Medicaid SPA'
$ws.Cells.Item(8, 4).Copy($ws.Cells.Item(73, 4))
$ws.Cells.Item(73, 5).Value = 'MD-25-9550'
$ws.Cells.Item(73, 6).Value = 'Submitted'
$ws.Cells.Item(8, 7).Copy($ws.Cells.Item(73, 7))

$ws.Cells.Item(74, 1).Value = 'SPA'
$ws.Cells.Item(74, 2).Value = 'MD'
$ws.Cells.Item(74, 3).Value = 'Medicaid SPA'
$ws.Cells.Item(8, 4).Copy($ws.Cells.Item(74, 4))
$ws.Cells.Item(74, 5).Value = 'MD-25-9551'
$ws.Cells.Item(74, 6).Value = 'Submitted'
$ws.Cells.Item(8, 7).Copy($ws.Cells.Item(74, 7))

$ws.Cells.Item(75, 1).Value = 'Waiver'
$ws.Cells.Item(75, 2).Value = 'MD'
$ws.Cells.Item(75, 3).Value = '1915(c)'
$ws.Cells.Item(75, 4).Value = 'Amendment'
$ws.Cells.Item(75, 5).Value = 'MD-2260.R00.69'
$ws.Cells.Item(75, 6).Value = 'Unsubmitted'
$ws.Cells.Item(75, 7).Value = 'MD-2260.R00.00'

$ws.Cells.Item(76, 1).Value = 'SPA'
$ws.Cells.Item(76, 2).Value = 'MD'
$ws.Cells.Item(76, 3).Value = 'Medicaid SPA'
$ws.Cells.Item(8, 4).Copy($ws.Cells.Item(76, 4))
$ws.Cells.Item(76, 5).Value = 'MD-25-9552'
$ws.Cells.Item(76, 6).Value = 'Under Review'
$ws.Cells.Item(8, 7).Copy($ws.Cells.Item(76, 7))

$ws.Cells.Item(77, 1).Value = 'Waiver'
$ws.Cells.Item(77, 2).Value = 'MD'
$ws.Cells.Item(77, 3).Value = '1915(b)'
$ws.Cells.Item(77, 4).Value = 'Initial'
$ws.Cells.Item(77, 5).Value = 'MD-2285.R00.00'
$ws.Cells.Item(77, 6).Value = 'Terminated'
$ws.Cells.Item(8, 7).Copy($ws.Cells.Item(77, 7))

$ws.Cells.Item(78, 1).Value = 'SPA'
$ws.Cells.Item(78, 2).Value = 'MD'
$ws.Cells.Item(78, 3).Value = 'Medicaid SPA'
$ws.Cells.Item(8, 4).Copy($ws.Cells.Item(78, 4))
$ws.Cells.Item(78, 5).Value = 'MD-25-9553'
$ws.Cells.Item(78, 6).Value = 'Withdrawn'
$ws.Cells.Item(8, 7).Copy($ws.Cells.Item(78, 7))
